$wb = $excel.ActiveWorkbook

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7642.0586
$ws.Range("I113").Value = 3446.2727
$ws.Range("J113").Value = 15334.333
$ws.Range("K113").Value = 3446.2727
$ws.Range("L113").Value = 15334.333
$ws.Range("M113").Value = -192.2727
$ws.Range("N113").Value = -21842.333

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 599365.0600000001
$ws.Range("I116").Value = 2003881.2
$ws.Range("J116").Value = 14150
$ws.Range("K116").Value = 2003881.2
$ws.Range("L116").Value = 14150
$ws.Range("M116").Value = -2000439.2
$ws.Range("N116").Value = -21034

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3796.5557
$ws.Range("I125").Value = 3366.4
$ws.Range("J125").Value = 3962
$ws.Range("K125").Value = 30297.6
$ws.Range("L125").Value = 35658
$ws.Range("M125").Value = -27837.6
$ws.Range("N125").Value = -40578

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 921.46
$ws.Range("J129").Value = 961.4681
$ws.Range("L129").Value = 2884.4043
$ws.Range("N129").Value = -12884.4043

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 26321166
$ws.Range("I132").Value = 31254978
$ws.Range("K132").Value = 93764934
$ws.Range("M132").Value = -93762404

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3497.2986
$ws.Range("I138").Value = 1925.1818
$ws.Range("J138").Value = 3806.1072
$ws.Range("K138").Value = 5775.5454
$ws.Range("L138").Value = 11418.3216
$ws.Range("M138").Value = -635.5454
$ws.Range("N138").Value = -21698.3216

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4758.0117
$ws.Range("I32").Value = 3114.1147
$ws.Range("J32").Value = 8936.25
$ws.Range("K32").Value = 3114.1147
$ws.Range("L32").Value = 8936.25
$ws.Range("M32").Value = -2827.1147
$ws.Range("N32").Value = -9510.25

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 795.63635
$ws.Range("I45").Value = 663
$ws.Range("K45").Value = 663
$ws.Range("M45").Value = -286

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2105.7144
$ws.Range("I61").Value = 1716.3636
$ws.Range("J61").Value = 3533.3333
$ws.Range("K61").Value = 1716.3636
$ws.Range("L61").Value = 3533.3333
$ws.Range("M61").Value = -1504.3636
$ws.Range("N61").Value = -3957.3333

# ARM row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 8429
$ws.Range("J123").Value = 8429
$ws.Range("L123").Value = 8429
$ws.Range("N123").Value = -18229

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2105.7144
$ws.Range("I136").Value = 1716.3636
$ws.Range("J136").Value = 3533.3333
$ws.Range("K136").Value = 5149.0908
$ws.Range("L136").Value = 10599.9999
$ws.Range("M136").Value = -2599.0908
$ws.Range("N136").Value = -15699.9999

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1664.7333
$ws.Range("I86").Value = 1574.6923
$ws.Range("J86").Value = 2250
$ws.Range("K86").Value = 1574.6923
$ws.Range("L86").Value = 2250
$ws.Range("M86").Value = -451.6922999999999
$ws.Range("N86").Value = -4496

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1664.7333
$ws.Range("I89").Value = 1574.6923
$ws.Range("J89").Value = 2250
$ws.Range("K89").Value = 7873.461499999999
$ws.Range("L89").Value = 11250
$ws.Range("M89").Value = -2257.461499999999
$ws.Range("N89").Value = -22482

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1800.25
$ws.Range("I107").Value = 1200.5
$ws.Range("K107").Value = 1200.5
$ws.Range("M107").Value = 719.5

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10529501
$ws.Range("I99").Value = 15385925
$ws.Range("J99").Value = 7249
$ws.Range("K99").Value = 15385925
$ws.Range("L99").Value = 7249
$ws.Range("M99").Value = -15384427
$ws.Range("N99").Value = -10245

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10529501
$ws.Range("I126").Value = 15385925
$ws.Range("J126").Value = 7249
$ws.Range("K126").Value = 46157775
$ws.Range("L126").Value = 21747
$ws.Range("M126").Value = -46155305
$ws.Range("N126").Value = -26687

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6259.7334
$ws.Range("I70").Value = 5967.76
$ws.Range("J70").Value = 6624.7
$ws.Range("K70").Value = 5967.76
$ws.Range("L70").Value = 6624.7
$ws.Range("M70").Value = -5697.76
$ws.Range("N70").Value = -7164.7

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6259.7334
$ws.Range("I73").Value = 5967.76
$ws.Range("J73").Value = 6624.7
$ws.Range("K73").Value = 5967.76
$ws.Range("L73").Value = 6624.7
$ws.Range("M73").Value = -5031.76
$ws.Range("N73").Value = -8496.700000000001

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2109.96
$ws.Range("I102").Value = 1548.7
$ws.Range("J102").Value = 2484.1333
$ws.Range("K102").Value = 1548.7
$ws.Range("L102").Value = 2484.1333
$ws.Range("M102").Value = 73.29999999999995
$ws.Range("N102").Value = -5728.1333

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2906.66
$ws.Range("I126").Value = 2921.2422
$ws.Range("J126").Value = 2214
$ws.Range("K126").Value = 8763.7266
$ws.Range("L126").Value = 6642
$ws.Range("M126").Value = -6293.7266
$ws.Range("N126").Value = -11582

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3855.5557
$ws.Range("I7").Value = 2333.4443
$ws.Range("J7").Value = 4616.6113
$ws.Range("K7").Value = 2333.4443
$ws.Range("L7").Value = 4616.6113
$ws.Range("M7").Value = -2221.4443
$ws.Range("N7").Value = -4840.6113

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5433.143
$ws.Range("I40").Value = 3649.0938
$ws.Range("K40").Value = 3649.0938
$ws.Range("M40").Value = -3513.0938

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5217.9375
$ws.Range("I122").Value = 3048.8
$ws.Range("K122").Value = 9146.400000000001
$ws.Range("M122").Value = -6696.400000000001

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3855.5557
$ws.Range("I126").Value = 2333.4443
$ws.Range("J126").Value = 4616.6113
$ws.Range("K126").Value = 7000.3329
$ws.Range("L126").Value = 13849.8339
$ws.Range("M126").Value = -4530.3329
$ws.Range("N126").Value = -18789.8339

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5235.7144
$ws.Range("I136").Value = 1716.6666
$ws.Range("K136").Value = 5149.9998
$ws.Range("M136").Value = -2599.9998

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 36853.8
$ws.Range("J123").Value = 36853.8
$ws.Range("L123").Value = 36853.8
$ws.Range("N123").Value = -46653.8

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2128.1592
$ws.Range("I126").Value = 1318.8846
$ws.Range("J126").Value = 3297.111
$ws.Range("K126").Value = 3956.6538
$ws.Range("L126").Value = 9891.332999999999
$ws.Range("M126").Value = -1486.6538
$ws.Range("N126").Value = -14831.333
